$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worker "KEVIN BARRIOS OSORIO" (CC 1143373237) rows move up to the top
# (rows 16-22) with periods in descending order (2211..2205), and the
# previously-first worker "ANDRESON EXMIT ZUÑIGA ARISTIZABAL" (CC 1201256725,
# period 1802) moves down to become the last row (23).

$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1143373237"
$ws.Range("D16").Value = "KEVIN BARRIOS OSORIO"
$ws.Range("E16").Value = "2211"
$ws.Range("F16").Value = 37800
$ws.Range("G16").Value = 1350000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1143373237"
$ws.Range("D17").Value = "KEVIN BARRIOS OSORIO"
$ws.Range("E17").Value = "2210"
$ws.Range("F17").Value = 54000
$ws.Range("G17").Value = 1350000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1143373237"
$ws.Range("D18").Value = "KEVIN BARRIOS OSORIO"
$ws.Range("E18").Value = "2209"
$ws.Range("F18").Value = 54000
$ws.Range("G18").Value = 1350000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1143373237"
$ws.Range("D19").Value = "KEVIN BARRIOS OSORIO"
$ws.Range("E19").Value = "2208"
$ws.Range("F19").Value = 54000
$ws.Range("G19").Value = 1350000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1143373237"
$ws.Range("D20").Value = "KEVIN BARRIOS OSORIO"
$ws.Range("E20").Value = "2207"
$ws.Range("F20").Value = 54000
$ws.Range("G20").Value = 1350000

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1143373237"
$ws.Range("D21").Value = "KEVIN BARRIOS OSORIO"
$ws.Range("E21").Value = "2206"
$ws.Range("F21").Value = 54000
$ws.Range("G21").Value = 1350000

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1143373237"
$ws.Range("D22").Value = "KEVIN BARRIOS OSORIO"
$ws.Range("E22").Value = "2205"
$ws.Range("F22").Value = 54000
$ws.Range("G22").Value = 1350000

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1201256725"
$ws.Range("D23").Value = "ANDRESON EXMIT ZUÑIGA ARISTIZABAL"
$ws.Range("E23").Value = "1802"
$ws.Range("F23").Value = 72000
$ws.Range("G23").Value = 1800000
